# getDEmand: refresh per-station demand start-times and realign the
# saved cursor/selection on every Station sheet (Station5 becomes the
# sheet that is on top / active when the workbook is saved).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Station1 - shift the 09:00 slot-end to 09:01 (row 2, col B)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Station1")
$ws1.Range("B2").Value = 0.3756944444444445

# ---------------------------------------------------------------
# Station2 - same 09:00 -> 09:01 shift
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Station2")
$ws2.Range("B2").Value = 0.3756944444444445

# ---------------------------------------------------------------
# Station3 - 09:00 -> 09:10 and 10:00 -> 10:10
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Station3")
$ws3.Range("B2").Value = 0.38194444444444442
$ws3.Range("B3").Value = 0.4236111111111111

# ---------------------------------------------------------------
# Station4 - 09:00 -> 09:01 and 10:00 -> 10:01
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Station4")
$ws4.Range("B2").Value = 0.3756944444444445
$ws4.Range("B3").Value = 0.41736111111111113

# ---------------------------------------------------------------
# Station5 - 10:00 -> 10:01 and 11:00 -> 11:01
# ---------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Station5")
$ws5.Range("B3").Value = 0.41736111111111113
$ws5.Range("B4").Value = 0.45902777777777781

# ---------------------------------------------------------------
# Re-point the saved selection on every sheet to match where the
# author left the cursor. Selecting on a sheet also makes it the
# active one, so Station5 (the sheet that should stay "on top")
# is selected last.
# ---------------------------------------------------------------
$ws1.Range("D13").Select()
$ws2.Range("D17").Select()
$ws3.Range("D15").Select()
$ws4.Range("D20").Select()
$ws5.Range("D17").Select()
$ws5.Activate()
